$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# The title currently holds three runs ("Lists", " ", "(continued)") whose
# concatenation already equals the desired final string, so a single
# direct assignment to "Lists (continued)" gets treated as a no-op diff.
# Route through an unrelated placeholder first so the engine actually
# rewrites the paragraph into one run, then set the real text.
$tr.Text = "placeholder"
$tr.Text = "Lists (continued)"
